$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in Count - Low (B), Count - Moderate (C), Count - High (D) for each category row
$data = @(
    @(2, 876, 3430, 672),
    @(3, 41, 300, 125),
    @(4, 2, 20, 4),
    @(5, 5, 14, 3),
    @(6, 28, 299, 72),
    @(7, 229, 1634, 801),
    @(8, 288, 1184, 246),
    @(9, 61, 327, 79),
    @(10, 958, 5487, 815),
    @(11, 231, 1060, 141),
    @(12, 40, 160, 12),
    @(13, 47, 456, 272),
    @(14, 150, 907, 374),
    @(15, 253, 1509, 471),
    @(16, 755, 6428, 2967),
    @(17, 2211, 11565, 4400)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Update the active selection to match the last-saved view state
$ws.Range("H17").Select()
